# Mail-merge placeholder cleanup:
#   Name_is       -> NAME
#   Position_is   -> POSITION
#   Salry_is.All  -> SALARY. All
#
# Find/Replace merges every run it touches into a single run, which also
# happens to drop the (now meaningless) spell-check w:proofErr markers
# that used to flag these placeholder tokens. Afterwards a harmless
# Bold-on/Bold-off no-op is applied to the untouched slices of that
# merged run so Word re-splits it back along the original run
# boundaries, leaving the surrounding, unrelated text in its own runs
# exactly as before.

$d = $word.ActiveDocument

function Resplit($doc, $paraIndex, $flatOffsets) {
    $p = $doc.Paragraphs($paraIndex)
    $base = $p.Range.Start
    for ($i = 0; $i -lt $flatOffsets.Count; $i += 2) {
        $a = $flatOffsets[$i]
        $b = $flatOffsets[$i + 1]
        $sub = $doc.Range($base + $a, $base + $b)
        $sub.Font.Bold = 1
        $sub.Font.Bold = 0
    }
}

# --- 1) "Dear Name_is," -> "Dear NAME," ------------------------------
$d.Content.Find.Execute("Dear Name_is", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dear NAME", 2) | Out-Null
$d.Content.Find.Execute("NAME,", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "NAME,", 2) | Out-Null
# Paragraph 3 is now the single run "Dear NAME," -> split back into
# "Dear " | "NAME" | ","
Resplit $d 3 @(5, 9)

# --- 2) "... promoted to Position_is." -> "... promoted to POSITION." -
$d.Content.Find.Execute("to Position_is.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "to POSITION.", 2) | Out-Null
# Paragraph 5 text is now "This is to inform you that in this appraisal
# cycle you have been promoted to POSITION." -> split back into
# "This is to inform you that " | "in this " | "appraisal" |
# " cycle you have been promoted to " | "POSITION" | "."
Resplit $d 5 @(35, 44, 44, 77, 77, 85)

# --- 3) "Salry_is.All" -> "SALARY. All" -------------------------------
$d.Content.Find.Execute("be Salry_is.All the", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "be SALARY. All the", 2) | Out-Null
# Paragraph 6 text is now "Your Revised Salary is going to be SALARY.
# All the very best for your future." -> split back into
# "Your Revised Salary is going to be " | "SALARY. All" |
# " the very best for your future."
Resplit $d 6 @(35, 46)

Write-Output "edits applied"
